$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Company Number) and I (SIC Codes) can contain purely-numeric
# looking text (e.g. "16473606", "64209") that must stay text, matching the
# original inlineStr string cells. Force text number format before writing
# so Excel doesn't silently coerce them to numeric cells.
$ws.Range("B3:B11").NumberFormat = "@"
$ws.Range("I3:I11").NumberFormat = "@"

# Row 3: SEVEN (HOLDCO) LIMITED
$ws.Range("A3").Value = "SEVEN (HOLDCO) LIMITED"
$ws.Range("B3").Value = "16473606"
$ws.Range("H3").Value = "Other"
$ws.Range("I3").Value = "64209"
$ws.Range("J3").Value = "Activities of other holding companies n.e.c."
$ws.Range("K3").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 4: GANDER INVESTMENTS LTD
$ws.Range("A4").Value = "GANDER INVESTMENTS LTD"
$ws.Range("B4").Value = "16473515"
$ws.Range("H4").Value = "Investments"
$ws.Range("I4").Value = "68100,68209"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""

# Row 5: INTERCONTINENTAL HOLDING COMPANY LIMITED
$ws.Range("A5").Value = "INTERCONTINENTAL HOLDING COMPANY LIMITED"
$ws.Range("B5").Value = "16473418"
$ws.Range("H5").Value = "Other"
$ws.Range("I5").Value = "64209"
$ws.Range("J5").Value = "Activities of other holding companies n.e.c."
$ws.Range("K5").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 6: TLJ INVESTMENT LTD
$ws.Range("A6").Value = "TLJ INVESTMENT LTD"
$ws.Range("B6").Value = "16473151"
$ws.Range("H6").Value = "Investments"
$ws.Range("I6").Value = "41100,55100,68100"
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""

# Row 7: GAUNT CAPITAL LTD
$ws.Range("A7").Value = "GAUNT CAPITAL LTD"
$ws.Range("B7").Value = "16473262"
$ws.Range("H7").Value = "Capital"
$ws.Range("I7").Value = "64209"
$ws.Range("J7").Value = "Activities of other holding companies n.e.c."
$ws.Range("K7").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 8: THE DISLEY GROUP LTD -- unchanged, left as-is

# Row 9: AJ INVESTMENT AND CONSULTANCY LTD
$ws.Range("A9").Value = "AJ INVESTMENT AND CONSULTANCY LTD"
$ws.Range("B9").Value = "16473328"
$ws.Range("I9").Value = "64306,70229"
$ws.Range("J9").Value = "Activities of real estate investment trusts"
$ws.Range("K9").Value = "UK-regulated REIT companies."

# Row 10: MARMIMI HOLDING LIMITED
$ws.Range("A10").Value = "MARMIMI HOLDING LIMITED"
$ws.Range("B10").Value = "16473234"

# Row 11: BRIDGEWICK PARTNERS LIMITED
$ws.Range("A11").Value = "BRIDGEWICK PARTNERS LIMITED"
$ws.Range("B11").Value = "16473142"
$ws.Range("H11").Value = "Partners"
$ws.Range("I11").Value = "64999"
$ws.Range("J11").Value = "Financial intermediation not elsewhere classified"
$ws.Range("K11").Value = "Catch-all credit-oriented SPVs for novel lending structures."
